# Update the build version string everywhere it appears in the workbook.
# Old: "mines - January 30 (built on January 30 2026 16.19.47 EST)"
# New: "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---
# A2: "Version: mines - January 30 (built on January 30 2026 16.19.47 EST)"
$wsAbout.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended Citation sentence embedding the version string.
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Bulga Coal Mine, Australia, M0017, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
# Column S ("build_version") for rows 2 through 37 holds the version string.
$usedRange = $wsData.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
